{"js": "// Replace each two-digit-multiplication expression in the document with\n// its new value, as described by the diff. Every occurrence is unique in\n// the document, so a simple search + full-text replace per pair is safe.\nconst replacements = [\n  [\"90\u00d735=\", \"11\u00d773=\"],\n  [\"33\u00d756=\", \"24\u00d742=\"],\n  [\"26\u00d743=\", \"20\u00d720=\"],\n  [\"38\u00d755=\", \"38\u00d781=\"],\n  [\"27\u00d736=\", \"57\u00d719=\"],\n  [\"42\u00d762=\", \"40\u00d727=\"],\n  [\"85\u00d762=\", \"96\u00d722=\"],\n  [\"96\u00d752=\", \"95\u00d712=\"],\n  [\"29\u00d733=\", \"84\u00d745=\"],\n  [\"68\u00d714=\", \"41\u00d762=\"],\n  [\"44\u00d735=\", \"38\u00d735=\"],\n  [\"19\u00d761=\", \"62\u00d761=\"],\n  [\"88\u00d779=\", \"57\u00d788=\"],\n  [\"71\u00d765=\", \"77\u00d775=\"],\n  [\"80\u00d713=\", \"24\u00d790=\"],\n  [\"71\u00d758=\", \"44\u00d748=\"],\n  [\"94\u00d719=\", \"97\u00d778=\"],\n  [\"38\u00d767=\", \"38\u00d724=\"],\n  [\"83\u00d773=\", \"43\u00d759=\"],\n  [\"37\u00d772=\", \"59\u00d785=\"],\n  [\"28\u00d736=\", \"46\u00d722=\"],\n  [\"87\u00d730=\", \"23\u00d799=\"],\n  [\"33\u00d773=\", \"92\u00d758=\"],\n  [\"26\u00d789=\", \"94\u00d797=\"],\n  [\"19\u00d727=\", \"14\u00d735=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit-multiplication expression in the document with\n# its new value, as described by the diff. Every occurrence is unique in\n# the document, so a Find/Replace pass per pair is safe and unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @{old=\"90\u00d735=\"; new=\"11\u00d773=\"},\n    @{old=\"33\u00d756=\"; new=\"24\u00d742=\"},\n    @{old=\"26\u00d743=\"; new=\"20\u00d720=\"},\n    @{old=\"38\u00d755=\"; new=\"38\u00d781=\"},\n    @{old=\"27\u00d736=\"; new=\"57\u00d719=\"},\n    @{old=\"42\u00d762=\"; new=\"40\u00d727=\"},\n    @{old=\"85\u00d762=\"; new=\"96\u00d722=\"},\n    @{old=\"96\u00d752=\"; new=\"95\u00d712=\"},\n    @{old=\"29\u00d733=\"; new=\"84\u00d745=\"},\n    @{old=\"68\u00d714=\"; new=\"41\u00d762=\"},\n    @{old=\"44\u00d735=\"; new=\"38\u00d735=\"},\n    @{old=\"19\u00d761=\"; new=\"62\u00d761=\"},\n    @{old=\"88\u00d779=\"; new=\"57\u00d788=\"},\n    @{old=\"71\u00d765=\"; new=\"77\u00d775=\"},\n    @{old=\"80\u00d713=\"; new=\"24\u00d790=\"},\n    @{old=\"71\u00d758=\"; new=\"44\u00d748=\"},\n    @{old=\"94\u00d719=\"; new=\"97\u00d778=\"},\n    @{old=\"38\u00d767=\"; new=\"38\u00d724=\"},\n    @{old=\"83\u00d773=\"; new=\"43\u00d759=\"},\n    @{old=\"37\u00d772=\"; new=\"59\u00d785=\"},\n    @{old=\"28\u00d736=\"; new=\"46\u00d722=\"},\n    @{old=\"87\u00d730=\"; new=\"23\u00d799=\"},\n    @{old=\"33\u00d773=\"; new=\"92\u00d758=\"},\n    @{old=\"26\u00d789=\"; new=\"94\u00d797=\"},\n    @{old=\"19\u00d727=\"; new=\"14\u00d735=\"}\n)\n\nforeach ($p in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $p.old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $p.new\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
